# Auto-generated edit script: update cryptos Price (D) / Volume(1h) (E) columns
# to match the "Updated cryptos list" GitHub Actions commit.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'43.593.46"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  -1.20%  "
$ws.Range("D3").Value = "'2.233.54"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +0.36%  "
$ws.Range("E4").Value = "  +0.08%  "
$ws.Range("D5").Value = "'269.25"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +3.31%  "
$ws.Range("D6").Value = "'94.19"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +13.42%  "
$ws.Range("D7").Value = "'0.622"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -1.18%  "
$ws.Range("E8").Value = "  +0.07%  "
$ws.Range("D9").Value = "'0.622"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +1.64%  "
$ws.Range("D10").Value = "'46.07"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +3.81%  "
$ws.Range("E11").Value = "  -1.23%  "
$ws.Range("D12").Value = "'8.23"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +16.29%  "
$ws.Range("E13").Value = "  +1.12%  "
$ws.Range("D14").Value = "'15.14"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +3.07%  "
$ws.Range("D15").Value = "'2.567.51"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +0.25%  "
$ws.Range("D16").Value = "'2.226.16"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +0.12%  "
$ws.Range("E17").Value = "  +1.48%  "
$ws.Range("D18").Value = "'43.562.67"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -1.01%  "
$ws.Range("D19").Value = "'0.0000104"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -1.20%  "
$ws.Range("E20").Value = "  -0.52%  "
$ws.Range("E21").Value = "  -1.74%  "
$ws.Range("E22").Value = "  -1.79%  "
$ws.Range("D23").Value = "'233.46"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -0.12%  "
$ws.Range("D24").Value = "'9.03"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -2.49%  "
$ws.Range("D26").Value = "'2.50"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +11.15%  "
$ws.Range("D27").Value = "'11.27"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +4.26%  "
$ws.Range("D28").Value = "'3.54"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +5.36%  "
$ws.Range("D29").Value = "'40.46"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -0.80%  "
$ws.Range("E30").Value = "  +1.96%  "
$ws.Range("D31").Value = "'173.12"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -0.18%  "
$ws.Range("E32").Value = "  +3.78%  "
$ws.Range("E33").Value = "  +0.72%  "
$ws.Range("D34").Value = "'5.47"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +1.89%  "
$ws.Range("E35").Value = "  -0.15%  "
$ws.Range("E36").Value = "  -4.85%  "
$ws.Range("E37").Value = "  -4.08%  "
$ws.Range("D38").Value = "'4.34"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -4.12%  "
$ws.Range("D39").Value = "'3.61"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +19.72%  "
$ws.Range("D40").Value = "'12.62"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -6.91%  "
$ws.Range("E41").Value = "  +2.72%  "
$ws.Range("D42").Value = "'0.219"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +7.72%  "
$ws.Range("D43").Value = "'63.28"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -1.17%  "
$ws.Range("E44").Value = "  -3.56%  "
$ws.Range("D45").Value = "'0.0989"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +0.13%  "
$ws.Range("D46").Value = "'8.38"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -0.05%  "
$ws.Range("D47").Value = "'100.78"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -2.86%  "
$ws.Range("D48").Value = "'1.15"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +1.95%  "
$ws.Range("E49").Value = "  +2.65%  "
$ws.Range("D50").Value = "'0.439"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -2.34%  "
$ws.Range("D51").Value = "'2.454.76"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +0.32%  "
